$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 376 (pushes the existing rows
# 376..430 down to 377..431, carrying their data/formatting with them).
$ws.Rows.Item(376).Insert()

# Populate the new row 376 with this week's data point (the newest
# "Ciboulette" - Vega Central Mapocho de Santiago record).
$ws.Cells.Item(376, 1).Value = 9
$ws.Cells.Item(376, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(376, 3).Value = "Metropolitana"
$ws.Cells.Item(376, 4).Value = 44776
$ws.Cells.Item(376, 5).Value = 13
$ws.Cells.Item(376, 6).Value = 100112039
$ws.Cells.Item(376, 7).Value = "Ciboulette"
$ws.Cells.Item(376, 8).Value = "Sin especificar"
$ws.Cells.Item(376, 9).Value = "Primera"
$ws.Cells.Item(376, 10).Value = 250
$ws.Cells.Item(376, 11).Value = 2000
$ws.Cells.Item(376, 12).Value = 2200
$ws.Cells.Item(376, 13).Value = 2100
$ws.Cells.Item(376, 14).Value = "$/docena de atados"
$ws.Cells.Item(376, 15).Value = "Región Metropolitana"
$ws.Cells.Item(376, 16).Value = 700
$ws.Cells.Item(376, 17).Value = 3
$ws.Cells.Item(376, 18).Value = "Hortaliza"
